# [ADD] Funcion para importar registro de costos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old data rows (A2:K3) so no stale values/styles bleed
#     through. The rest of row 1 (A1:I1) already carries the
#     bold/border/center style ("s=1") and is only having its text
#     replaced below, so leave its formatting untouched. ---
$ws.Range("A2:K3").Clear()

# --- Columns J and K are no longer part of the table; drop them
#     entirely (not just their contents) so the <cols> definition ends
#     at column I, matching the new 9-column layout. ---
$ws.Columns("J:K").Delete()

# --- Column widths: only columns A-I remain, all width 30.
#     NB: COM's ColumnWidth setter round-trips through a pixel-based MDW
#     rounding (like real Excel), so asking for "30" literally lands the
#     stored OOXML width at 30.8333. 29.125 is the character-width input
#     that this rounding maps back down to an exact stored width of 30. ---
$cols = @("A","B","C","D","E","F","G","H","I")
foreach ($col in $cols) {
    $ws.Range("$col`1").ColumnWidth = 29.125
}

# --- Header row (row 1) text, keeps its existing bold/border/center style ---
$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Comprobante"
$ws.Range("C1").Value = "Cliente / (Producto / Servicio)"
$ws.Range("D1").Value = "Estado"
$ws.Range("E1").Value = "Vendedor"
$ws.Range("F1").Value = "Guias"
$ws.Range("G1").Value = "Moneda"
$ws.Range("H1").Value = "Monto de la boleta"
$ws.Range("I1").Value = "Monto (S/)"

# --- Data rows: each row is a plain array of text values (col -> value); "" means "leave blank" ---
$data = @(
    @("2023-04-28", "B001-0028", "MARIA LUZ PEREZ DE RAMOS", "Aceptado", "USR-16", "['T001-0617']", "SOLES", "512.18", "512.18"),
    @("2023-04-26", "B001-0027", "METALPROTEC S.A.C ", "", "USR-16", "['T001-0613']", "SOLES", "1042.15", "1042.15"),
    @("2023-04-26", "B001-0026", "METALPROTEC S.A.C ", "", "USR-16", "['T001-0611']", "SOLES", "733.46", "733.46"),
    @("2023-04-24", "B001-0025", "ELIAS MARCHENA MARCHENA", "Aceptado", "USR-16", "['T001-0603']", "SOLES", "258.21", "258.21"),
    @("2023-04-19", "B001-0024", "ELIAS MARCHENA MARCHENA", "Aceptado", "USR-16", "['T001-0593']", "SOLES", "259.11", "259.11"),
    @("2023-04-10", "B001-0023", "SEGUNDO CARCAMO NIZAMA", "Aceptado", "USR-16", "['T001-0573']", "SOLES", "139.04", "139.04"),
    @("2023-04-01", "B001-0022", "EUSTACIO VIDAL CESPEDES", "Aceptado", "USR-16", "['T001-0561']", "SOLES", "154.77", "154.77")
)

$rowIndex = 2
foreach ($rowData in $data) {
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $value = $rowData[$c]
        if ($value -ne "") {
            $cell = $ws.Cells.Item($rowIndex, $c + 1)
            # Force text storage so date-like / numeric-like strings aren't
            # reinterpreted as dates or numbers (matches source data which
            # is plain text throughout).
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.ClearFormats()
        }
    }
    $rowIndex++
}

# --- Total row (row 9) ---
$totalLabelCell = $ws.Cells.Item(9, 8)
$totalLabelCell.NumberFormat = "@"
$totalLabelCell.Value = "Monto Total"
$totalLabelCell.ClearFormats()

$totalValueCell = $ws.Cells.Item(9, 9)
$totalValueCell.NumberFormat = "@"
$totalValueCell.Value = "1323.31"
$totalValueCell.ClearFormats()
